# Edit: rename input/variable headers and update category & condition
# (tek shares) labels on the s_curve sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) renames -------------------------------------------
# A1: building_category (unchanged text, was already building_category)
# B1: renovation_type -> condition
# C1: earliest_renovation_age -> earliest_age_for_measure
# D1: average_age -> average_age_for_measure
# E1: rush_period_years (unchanged)
# F1: last_renovation_age -> last_age_for_measure
# G1: rush_share (unchanged)
# H1: never_share (unchanged)
$ws.Range("A1").Value = "building_category"
$ws.Range("B1").Value = "condition"
$ws.Range("C1").Value = "earliest_age_for_measure"
$ws.Range("D1").Value = "average_age_for_measure"
$ws.Range("F1").Value = "last_age_for_measure"

# --- Rename building category (column A) and condition/measure (column B) values ---
$nameMapA = @{
    "Apartment"      = "Apartment block"
    "SmallHouse"     = "House"
    "Shop"           = "Retail"
    "StorageRepairs" = "Storage repairs"
}
$nameMapB = @{
    "SmallMeasures"  = "Small measure"
    "Rehabilitation" = "Renovation"
}

for ($r = 2; $r -le 40; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($nameMapA.ContainsKey($aVal)) {
        $aCell.Value = $nameMapA[$aVal]
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($nameMapB.ContainsKey($bVal)) {
        $bCell.Value = $nameMapB[$bVal]
    }
}

# --- Column width adjustments (columns widened to fit the new, longer header text) ---
$ws.Columns.Item(3).ColumnWidth = 23.6666666666667
$ws.Columns.Item(4).ColumnWidth = 24
$ws.Columns.Item(6).ColumnWidth = 20
